$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove stale component-forecast values from rows 2-4 (naive forecaster bug fix)
$ws.Range("C2").ClearContents()
$ws.Range("E2").ClearContents()
$ws.Range("C3").ClearContents()
$ws.Range("C4").ClearContents()

# Refresh recomputed naive-forecast values (presentation state 11.02 recalculation)
$ws.Range("E4").Value = 3.386383090739975
$ws.Range("C5").Value = 0.8787496612563173
$ws.Range("E5").Value = 1.013823151053028
$ws.Range("C6").Value = 2.533533936850585
$ws.Range("E6").Value = 0.9842934829757288
$ws.Range("E7").Value = 4.356912452939454
$ws.Range("E8").Value = 3.612753212925446
$ws.Range("C9").Value = 2.431458940166964
$ws.Range("E9").Value = 5.259925231829876
$ws.Range("C10").Value = 1.21254482274098
$ws.Range("E10").Value = 2.158838189283219
$ws.Range("C11").Value = 1.447930496829564
$ws.Range("E12").Value = 1.194058515117313
$ws.Range("E13").Value = -0.563208905821222
$ws.Range("C14").Value = 0.4712609263772816
$ws.Range("E14").Value = 1.409662779709797
$ws.Range("C15").Value = 0.6742451383204839
$ws.Range("E15").Value = 1.713290556413583
$ws.Range("E18").Value = 4.595879021798344
$ws.Range("C20").Value = 4.109890522944326
$ws.Range("C21").Value = 1.715791310593251
$ws.Range("E21").Value = 1.687339605296501
$ws.Range("E22").Value = 0.02883756256673031
$ws.Range("C23").Value = 1.862609889357336
$ws.Range("E24").Value = -2.079848588862143
$ws.Range("E25").Value = 2.147322685428366
$ws.Range("E26").Value = 0.9262553939923146
$ws.Range("C28").Value = 0.893498267486792
$ws.Range("E28").Value = -1.194610791899986
$ws.Range("E30").Value = 2.928189816005689
$ws.Range("C31").Value = 2.306826470345391
$ws.Range("E31").Value = 1.40519946540949
$ws.Range("E32").Value = 0.8024032016000104
$ws.Range("C33").Value = 3.265677646667942
$ws.Range("E33").Value = 5.715169758465
$ws.Range("E34").Value = 3.828814763561783
$ws.Range("C35").Value = 0.2974381310041352
$ws.Range("E35").Value = -2.009776081564663
$ws.Range("E36").Value = 0.8023688159249032
$ws.Range("E37").Value = 6.778609849419737
$ws.Range("C38").Value = 2.777797690741446
$ws.Range("C39").Value = 2.475264839201419
$ws.Range("C40").Value = 0.06579575777907465
$ws.Range("E40").Value = 0.2740865344839749
$ws.Range("C43").Value = 1.076435582022328
$ws.Range("E43").Value = 2.51686114938241
$ws.Range("C45").Value = -1.650648527511434
$ws.Range("E45").Value = -0.2006752520846145
$ws.Range("C46").Value = -1.432689847121826
$ws.Range("E46").Value = 0.1752798163574321
$ws.Range("E47").Value = 1.396505962682837
$ws.Range("E48").Value = -1.696610696428313
$ws.Range("C50").Value = 2.033479419175155
$ws.Range("C51").Value = 3.147579643557918
$ws.Range("E52").Value = -1.362365718491854
$ws.Range("C53").Value = 2.581716327283523
